$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36.. down by one.
$ws.Rows.Item(36).Insert()

# Fill the new row 36 with data (same constant columns as the rest of the table,
# plus the specific varying values for this entry).
$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(36, 3).Value = "Ñuble"
$ws.Cells.Item(36, 4).Value = 45002
$ws.Cells.Item(36, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 5).Value = 16
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100108
$ws.Cells.Item(36, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(36, 9).Value = 100108002
$ws.Cells.Item(36, 10).Value = "Mango"
$ws.Cells.Item(36, 11).Value = "Sin especificar"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 40
$ws.Cells.Item(36, 14).Value = 7500
$ws.Cells.Item(36, 15).Value = 8000
$ws.Cells.Item(36, 16).Value = 7750
$ws.Cells.Item(36, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(36, 18).Value = "Perú"
$ws.Cells.Item(36, 19).Value = 1938
$ws.Cells.Item(36, 20).Value = 4
